$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.253.79'
$ws.Range("E2").Value = '  -2.88%  '

$ws.Range("D3").Value = '2.991.73'
$ws.Range("E3").Value = '  -3.60%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.64'
$ws.Range("E5").Value = '  -2.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.59'
$ws.Range("E6").Value = '  -7.27%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.521'
$ws.Range("E8").Value = '  -3.30%  '

$ws.Range("D9").Value = '2.991.84'
$ws.Range("E9").Value = '  -3.55%  '

$ws.Range("E10").Value = '  -7.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.440'
$ws.Range("E12").Value = '  -2.63%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000227'
$ws.Range("E13").Value = '  -5.22%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.51'
$ws.Range("E14").Value = '  -6.76%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.123'
$ws.Range("E15").Value = '  +1.98%  '

$ws.Range("D16").Value = '3.485.24'
$ws.Range("E16").Value = '  -3.46%  '

$ws.Range("D17").Value = '62.267.08'
$ws.Range("E17").Value = '  -2.65%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.98'
$ws.Range("E18").Value = '  -3.31%  '

$ws.Range("D19").Value = '2.996.15'
$ws.Range("E19").Value = '  -3.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '456.53'
$ws.Range("E20").Value = '  -5.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.81'
$ws.Range("E21").Value = '  -4.79%  '

$ws.Range("E22").Value = '  -5.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.32'
$ws.Range("E23").Value = '  -3.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.00'
$ws.Range("E24").Value = '  -1.96%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.25'
$ws.Range("E25").Value = '  -9.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.17'
$ws.Range("E26").Value = '  -5.91%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.08'
$ws.Range("E27").Value = '  -6.70%  '

$ws.Range("E28").Value = '  -0.07%  '

$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.61'
$ws.Range("E30").Value = '  -3.50%  '

$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.12'
$ws.Range("E31").Value = '  -6.99%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.07'
$ws.Range("E32").Value = '  -6.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.78'
$ws.Range("E33").Value = '  -1.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.107'
$ws.Range("E34").Value = '  -5.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.02'
$ws.Range("E35").Value = '  -5.82%  '

$ws.Range("D36").Value = '0.0₃0783'
$ws.Range("E36").Value = '  -7.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.71'
$ws.Range("E37").Value = '  -5.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.10'
$ws.Range("E38").Value = '  -6.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.04'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.92'
$ws.Range("E40").Value = '  -3.47%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.92'
$ws.Range("E41").Value = '  -11.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '405.28'
$ws.Range("E42").Value = '  -9.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.114'
$ws.Range("E43").Value = '  +0.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.273'
$ws.Range("E44").Value = '  -6.49%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '39.17'
$ws.Range("E45").Value = '  -2.68%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.755.11'
$ws.Range("E46").Value = '  -3.10%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0349'
$ws.Range("E47").Value = '  -4.50%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.92'
$ws.Range("E48").Value = '  -3.28%  '

$ws.Range("E50").Value = '  -2.38%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.59'
$ws.Range("E51").Value = '  -9.72%  '
